$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 1154, pushing existing rows 1154:1239 down to 1156:1241.
$ws.Range("A1154:A1155").EntireRow.Insert()

# Populate new row 1154 (Betarraga, Primera, week of 44931)
$ws.Cells.Item(1154, 1).Value = 6
$ws.Cells.Item(1154, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(1154, 3).Value = "Metropolitana"
$ws.Cells.Item(1154, 4).Value = 44931
$ws.Cells.Item(1154, 5).Value = 13
$ws.Cells.Item(1154, 6).Value = 100114014
$ws.Cells.Item(1154, 7).Value = "Betarraga"
$ws.Cells.Item(1154, 8).Value = "Sin especificar"
$ws.Cells.Item(1154, 9).Value = "Primera"
$ws.Cells.Item(1154, 10).Value = 60000
$ws.Cells.Item(1154, 11).Value = 80
$ws.Cells.Item(1154, 12).Value = 90
$ws.Cells.Item(1154, 13).Value = 84
$ws.Cells.Item(1154, 14).Value = "`$/unidad"
$ws.Cells.Item(1154, 15).Value = "Región Metropolitana"
$ws.Cells.Item(1154, 16).Value = 84
$ws.Cells.Item(1154, 17).Value = 1
$ws.Cells.Item(1154, 18).Value = "Hortaliza"

# Populate new row 1155 (Betarraga, Segunda, week of 44931)
$ws.Cells.Item(1155, 1).Value = 6
$ws.Cells.Item(1155, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(1155, 3).Value = "Metropolitana"
$ws.Cells.Item(1155, 4).Value = 44931
$ws.Cells.Item(1155, 5).Value = 13
$ws.Cells.Item(1155, 6).Value = 100114014
$ws.Cells.Item(1155, 7).Value = "Betarraga"
$ws.Cells.Item(1155, 8).Value = "Sin especificar"
$ws.Cells.Item(1155, 9).Value = "Segunda"
$ws.Cells.Item(1155, 10).Value = 47000
$ws.Cells.Item(1155, 11).Value = 55
$ws.Cells.Item(1155, 12).Value = 60
$ws.Cells.Item(1155, 13).Value = 57
$ws.Cells.Item(1155, 14).Value = "`$/unidad"
$ws.Cells.Item(1155, 15).Value = "Región Metropolitana"
$ws.Cells.Item(1155, 16).Value = 57
$ws.Cells.Item(1155, 17).Value = 1
$ws.Cells.Item(1155, 18).Value = "Hortaliza"

# Apply the same numeric date format used by the rest of column D to the new D cells.
$ws.Cells.Item(1154, 4).NumberFormat = $ws.Cells.Item(1156, 4).NumberFormat
$ws.Cells.Item(1155, 4).NumberFormat = $ws.Cells.Item(1156, 4).NumberFormat

"done"
